$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark from its current location
# (end of the paragraph that ends with "...העכבר.") ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Step 2: append two new paragraphs after the last (empty) paragraph ---
# Word's terminal body paragraph can't be targeted with InsertXML directly
# (it gets merged into that paragraph), so use InsertParagraphAfter to make
# room, then InsertXML to overwrite that new paragraph's contents cleanly
# (no stray empty run).
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$dup = $r.Duplicate
$dup.Collapse(0)
$dup.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Last
$r1 = $newPara1.Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:bidi/></w:pPr></w:p>'
$r1.InsertXML($xml1)

$r1b = $d.Paragraphs.Last.Range
$dup2 = $r1b.Duplicate
$dup2.Collapse(0)
$dup2.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Last
$r2 = $newPara2.Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:bidi/><w:rPr><w:lang w:val="en-IL"/></w:rPr></w:pPr></w:p>'
$r2.InsertXML($xml2)

# --- Step 3: re-add the "_GoBack" bookmark, collapsed, on the new last paragraph ---
$finalPara = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $finalPara.Range)
